# Update "想去人数" (number of people interested) values in the
# 杭州-漫展信息 workbook across the four worksheets.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 4592
$ws.Cells.Item(3, 6).Value = 447
$ws.Cells.Item(4, 6).Value = 3705
$ws.Cells.Item(9, 6).Value = 389
$ws.Cells.Item(10, 6).Value = 2628
$ws.Cells.Item(11, 6).Value = 1298
$ws.Cells.Item(14, 6).Value = 289
$ws.Cells.Item(15, 6).Value = 29
$ws.Cells.Item(17, 6).Value = 268
$ws.Cells.Item(18, 6).Value = 74
$ws.Cells.Item(19, 6).Value = 10860
$ws.Cells.Item(20, 6).Value = 6235
$ws.Cells.Item(23, 6).Value = 399
$ws.Cells.Item(24, 6).Value = 236
$ws.Cells.Item(27, 6).Value = 862
$ws.Cells.Item(29, 6).Value = 210
$ws.Cells.Item(30, 6).Value = 874
$ws.Cells.Item(31, 6).Value = 3585
$ws.Cells.Item(33, 6).Value = 975
$ws.Cells.Item(36, 6).Value = 289
$ws.Cells.Item(37, 6).Value = 254
$ws.Cells.Item(38, 6).Value = 266
$ws.Cells.Item(39, 6).Value = 4899
$ws.Cells.Item(41, 6).Value = 1176
$ws.Cells.Item(43, 6).Value = 231
$ws.Cells.Item(44, 6).Value = 135
$ws.Cells.Item(45, 6).Value = 507

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(13, 6).Value = 3629
$ws.Cells.Item(21, 6).Value = 77

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 8903
$ws.Cells.Item(3, 6).Value = 458
$ws.Cells.Item(4, 6).Value = 1715

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 458
$ws.Cells.Item(3, 6).Value = 1715
$ws.Cells.Item(4, 6).Value = 4592
$ws.Cells.Item(5, 6).Value = 3705
$ws.Cells.Item(9, 6).Value = 389
$ws.Cells.Item(10, 6).Value = 2628
$ws.Cells.Item(14, 6).Value = 1298
$ws.Cells.Item(16, 6).Value = 289
$ws.Cells.Item(17, 6).Value = 29
$ws.Cells.Item(19, 6).Value = 268
$ws.Cells.Item(20, 6).Value = 10860
$ws.Cells.Item(21, 6).Value = 3629
$ws.Cells.Item(25, 6).Value = 399
$ws.Cells.Item(26, 6).Value = 236
$ws.Cells.Item(30, 6).Value = 210
$ws.Cells.Item(31, 6).Value = 874
$ws.Cells.Item(32, 6).Value = 3585
$ws.Cells.Item(34, 6).Value = 975
$ws.Cells.Item(36, 6).Value = 289
$ws.Cells.Item(37, 6).Value = 254
$ws.Cells.Item(40, 6).Value = 266
$ws.Cells.Item(41, 6).Value = 4899
$ws.Cells.Item(43, 6).Value = 1176
$ws.Cells.Item(45, 6).Value = 135
$ws.Cells.Item(46, 6).Value = 507
$ws.Cells.Item(47, 6).Value = 77
